# Fruta / hortaliza, semanal
# Insert a new weekly price-report block (4 rows) at the top of the
# "Piña" data table (rows 1086-1165), pushing the existing data down
# by 4 rows and updating the dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows before the current first data row of the block (1086).
# This shifts the existing rows 1086:1165 down to 1090:1169.
$ws.Rows("1086:1089").Insert()

# Constant values shared by every row in this subset (read off the row that
# used to be 1086, now shifted to 1090).
$mercadoId = $ws.Range("A1090").Value2
$mercado   = $ws.Range("B1090").Value2
$region    = $ws.Range("C1090").Value2
$codreg    = $ws.Range("E1090").Value2
$tipo      = $ws.Range("F1090").Value2
$productoId = $ws.Range("G1090").Value2
$producto  = $ws.Range("H1090").Value2
$categoriaId = $ws.Range("I1090").Value2
$categoria = $ws.Range("J1090").Value2
$variedad  = $ws.Range("K1090").Value2
$origen    = $ws.Range("R1090").Value2
$volumen   = $ws.Range("M1090").Value2

# New weekly block values (row 1086-1089), date 44826 (2022-09-22).
$fecha = 44826

$calidades = @("Especial", "Primera", "Segunda", "Tercera")
$precioMin = 21000
$precioMax = 22000
$precioProm = 21500
$unidades  = @("$/caja 10 unidades", "$/caja 12 unidades", "$/caja 14 unidades", "$/caja 16 unidades")
$precioKg  = @(2150, 1792, 1536, 1344)
$kgUnidad  = @(10, 12, 14, 16)

for ($i = 0; $i -lt 4; $i++) {
    $r = 1086 + $i

    $ws.Cells.Item($r, 1).Value2 = $mercadoId
    $ws.Cells.Item($r, 2).Value2 = $mercado
    $ws.Cells.Item($r, 3).Value2 = $region
    $ws.Range("D$r").Value2 = $fecha
    $ws.Cells.Item($r, 5).Value2 = $codreg
    $ws.Cells.Item($r, 6).Value2 = $tipo
    $ws.Cells.Item($r, 7).Value2 = $productoId
    $ws.Cells.Item($r, 8).Value2 = $producto
    $ws.Cells.Item($r, 9).Value2 = $categoriaId
    $ws.Cells.Item($r, 10).Value2 = $categoria
    $ws.Cells.Item($r, 11).Value2 = $variedad
    $ws.Cells.Item($r, 12).Value2 = $calidades[$i]
    $ws.Cells.Item($r, 13).Value2 = $volumen
    $ws.Cells.Item($r, 14).Value2 = $precioMin
    $ws.Cells.Item($r, 15).Value2 = $precioMax
    $ws.Cells.Item($r, 16).Value2 = $precioProm
    $ws.Cells.Item($r, 17).Value2 = $unidades[$i]
    $ws.Cells.Item($r, 18).Value2 = $origen
    $ws.Cells.Item($r, 19).Value2 = $precioKg[$i]
    $ws.Cells.Item($r, 20).Value2 = $kgUnidad[$i]
}
